$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels in row 1: "_A"/"_B" suffixes become randomized "_1"/"_2"
# (commit message: "randomized left and right order").
$ws.Range("A1").Value = "efficiency_1"
$ws.Range("G1").Value = "num_doses_2"
$ws.Range("B1").Value = "num_doses_1"
$ws.Range("C1").Value = "price_1"
$ws.Range("D1").Value = "num_X_1"
$ws.Range("E1").Value = "image_title_1"
$ws.Range("F1").Value = "efficiency_2"
$ws.Range("H1").Value = "price_2"
$ws.Range("I1").Value = "num_X_2"
$ws.Range("J1").Value = "image_title_2"

# Update the active cell selection to F2
$ws.Range("F2").Select()
